# Apply FHIR IG publisher regeneration changes to the "tddui-human-name"
# StructureDefinition summary workbook.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item(1)

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: regenerated timestamp
$meta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# Base Definition: pin the referenced profile to version 2.1.0
$meta.Range("B18").Value = "https://hl7.fr/ig/fhir/core/StructureDefinition/fr-core-human-name|2.1.0"

# ---- Elements sheet ---------------------------------------------------
$elem = $wb.Worksheets.Item(2)

# Type(s) column for the assemblyOrder extension row now pins the
# extension's version.
$elem.Range("K5").Value = "Extension {humanname-assembly-order|5.2.0}`n"

# Binding Value Set column: pin referenced value sets/code systems to a
# specific version.
$elem.Range("Z6").Value = "http://hl7.org/fhir/ValueSet/name-use|4.0.1"
$elem.Range("Z10").Value = "https://mos.esante.gouv.fr/NOS/JDV_J245-Civilite-CISIS/FHIR/JDV-J245-Civilite-CISIS|20230331120000"
$elem.Range("Z11").Value = "https://mos.esante.gouv.fr/NOS/JDV_J79-CiviliteExercice-RASS/FHIR/JDV-J79-CiviliteExercice-RASS|20200424120000"

# Column widths grew to fit the newly lengthened content in columns K
# (Type(s)) and Z (Binding Value Set). (Values chosen are the closest
# this engine's column-width quantization can reach to the authored
# widths of 38.39453125 / 91.87890625.)
$elem.Columns.Item(11).ColumnWidth = 37.42
$elem.Columns.Item(26).ColumnWidth = 90.92
